$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 210.85715
$ws.Range("I2").Value = 114.333336
$ws.Range("J2").Value = 790
$ws.Range("K2").Value = 114.333336
$ws.Range("L2").Value = 790
$ws.Range("M2").Value = -1.333336000000003
$ws.Range("N2").Value = -1016
# row 38
$ws.Range("H38").Value = 2075.2
$ws.Range("J38").Value = 3999
$ws.Range("L38").Value = 11997
$ws.Range("N38").Value = -12741
# row 40
$ws.Range("H40").Value = 600
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# row 43
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931
# row 70
$ws.Range("H70").Value = 2687.25
$ws.Range("J70").Value = 2687.25
$ws.Range("L70").Value = 8061.75
$ws.Range("N70").Value = -8601.75
# row 73
$ws.Range("H73").Value = 2687.25
$ws.Range("J73").Value = 2687.25
$ws.Range("L73").Value = 8061.75
$ws.Range("N73").Value = -9933.75
# row 80
$ws.Range("H80").Value = 789.95
$ws.Range("I80").Value = 735.2
$ws.Range("J80").Value = 844.7
$ws.Range("K80").Value = 2205.6
$ws.Range("L80").Value = 2534.1
$ws.Range("M80").Value = -1207.6
$ws.Range("N80").Value = -4530.1
# row 83
$ws.Range("H83").Value = 789.95
$ws.Range("I83").Value = 735.2
$ws.Range("J83").Value = 844.7
$ws.Range("K83").Value = 6616.8
$ws.Range("L83").Value = 7602.3
$ws.Range("M83").Value = -1624.8
$ws.Range("N83").Value = -17586.3
# row 98
$ws.Range("H98").Value = 3563.6
$ws.Range("J98").Value = 3318.6
$ws.Range("L98").Value = 3318.6
$ws.Range("N98").Value = -6314.6
# row 105
$ws.Range("H105").Value = 28000
$ws.Range("J105").Value = 28000
$ws.Range("L105").Value = 28000
$ws.Range("N105").Value = -34988
# row 122
$ws.Range("H122").Value = 3563.6
$ws.Range("J122").Value = 3318.6
$ws.Range("L122").Value = 9955.799999999999
$ws.Range("N122").Value = -14855.8
# row 127
$ws.Range("H127").Value = 501974.1
$ws.Range("I127").Value = 716627.3
$ws.Range("J127").Value = 1116.6666
$ws.Range("K127").Value = 2149881.9
$ws.Range("L127").Value = 3349.9998
$ws.Range("M127").Value = -2144921.9
$ws.Range("N127").Value = -13269.9998
# row 137
$ws.Range("H137").Value = 94220.25
$ws.Range("I137").Value = 1885.4
$ws.Range("J137").Value = 248111.67
$ws.Range("K137").Value = 5656.200000000001
$ws.Range("L137").Value = 744335.01
$ws.Range("M137").Value = -3106.200000000001
$ws.Range("N137").Value = -749435.01

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 1464.95
$ws.Range("I45").Value = 1293.4667
$ws.Range("J45").Value = 1979.4
$ws.Range("K45").Value = 1293.4667
$ws.Range("L45").Value = 1979.4
$ws.Range("M45").Value = -916.4666999999999
$ws.Range("N45").Value = -2733.4
# row 61
$ws.Range("H61").Value = 1662595.2
$ws.Range("I61").Value = 4502.067
$ws.Range("J61").Value = 7880444.5
$ws.Range("K61").Value = 4502.067
$ws.Range("L61").Value = 7880444.5
$ws.Range("M61").Value = -4290.067
$ws.Range("N61").Value = -7880868.5
# row 101
$ws.Range("H101").Value = 10000000
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# row 136
$ws.Range("H136").Value = 1662595.2
$ws.Range("I136").Value = 4502.067
$ws.Range("J136").Value = 7880444.5
$ws.Range("K136").Value = 13506.201
$ws.Range("L136").Value = 23641333.5
$ws.Range("M136").Value = -10956.201
$ws.Range("N136").Value = -23646433.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 86
$ws.Range("H86").Value = 9401.5
$ws.Range("I86").Value = 9939.666999999999
$ws.Range("K86").Value = 9939.666999999999
$ws.Range("M86").Value = -8816.666999999999
# row 89
$ws.Range("H89").Value = 9401.5
$ws.Range("I89").Value = 9939.666999999999
$ws.Range("K89").Value = 49698.335
$ws.Range("M89").Value = -44082.335
# row 122
$ws.Range("H122").Value = 1792.2727
$ws.Range("I122").Value = 1792.2727
$ws.Range("K122").Value = 5376.8181
$ws.Range("M122").Value = -2926.8181
# row 132
$ws.Range("H132").Value = 48130196
$ws.Range("I132").Value = 2165.2778
$ws.Range("K132").Value = 6495.8334
$ws.Range("M132").Value = -3965.8334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 673.4
$ws.Range("I5").Value = 674.2778
$ws.Range("J5").Value = 665.5
$ws.Range("K5").Value = 2022.8334
$ws.Range("L5").Value = 1996.5
$ws.Range("M5").Value = -1910.8334
$ws.Range("N5").Value = -2220.5
# row 49
$ws.Range("H49").Value = 1050
$ws.Range("J49").Value = 1050
$ws.Range("L49").Value = 3150
$ws.Range("N49").Value = -3462
# row 68
$ws.Range("H68").Value = 1136.4117
$ws.Range("I68").Value = 555
$ws.Range("J68").Value = 1315.3077
$ws.Range("K68").Value = 1665
$ws.Range("L68").Value = 3945.9231
$ws.Range("M68").Value = -854
$ws.Range("N68").Value = -5567.9231
# row 69
$ws.Range("H69").Value = 1798.5
$ws.Range("J69").Value = 2900
$ws.Range("L69").Value = 8700
$ws.Range("N69").Value = -10322
# row 71
$ws.Range("H71").Value = 1136.4117
$ws.Range("I71").Value = 555
$ws.Range("J71").Value = 1315.3077
$ws.Range("K71").Value = 4995
$ws.Range("L71").Value = 11837.7693
$ws.Range("M71").Value = -939
$ws.Range("N71").Value = -19949.7693
# row 72
$ws.Range("H72").Value = 1798.5
$ws.Range("J72").Value = 2900
$ws.Range("L72").Value = 26100
$ws.Range("N72").Value = -34212
# row 122
$ws.Range("H122").Value = 11327352
$ws.Range("I122").Value = 20764612
$ws.Range("J122").Value = 2833817
$ws.Range("K122").Value = 186881508
$ws.Range("L122").Value = 25504353
$ws.Range("M122").Value = -186879058
$ws.Range("N122").Value = -25509253
# row 135
$ws.Range("H135").Value = 673.4
$ws.Range("I135").Value = 674.2778
$ws.Range("J135").Value = 665.5
$ws.Range("K135").Value = 6068.500199999999
$ws.Range("L135").Value = 5989.5
$ws.Range("M135").Value = -3533.500199999999
$ws.Range("N135").Value = -11059.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 52
$ws.Range("H52").Value = 2276.24
$ws.Range("J52").Value = 2839.7896
$ws.Range("L52").Value = 2839.7896
$ws.Range("N52").Value = -3357.7896
# row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 4107.6665
$ws.Range("I22").Value = 916.6667
$ws.Range("J22").Value = 4905.4165
$ws.Range("K22").Value = 916.6667
$ws.Range("L22").Value = 4905.4165
$ws.Range("M22").Value = -621.6667
$ws.Range("N22").Value = -5495.4165
# row 27
$ws.Range("H27").Value = 4107.6665
$ws.Range("I27").Value = 916.6667
$ws.Range("J27").Value = 4905.4165
$ws.Range("K27").Value = 916.6667
$ws.Range("L27").Value = 4905.4165
$ws.Range("M27").Value = -809.6667
$ws.Range("N27").Value = -5119.4165
# row 55
$ws.Range("H55").Value = 1880.3334
$ws.Range("I55").Value = 1675.6
$ws.Range("J55").Value = 2136.25
$ws.Range("K55").Value = 1675.6
$ws.Range("L55").Value = 2136.25
$ws.Range("M55").Value = -1502.6
$ws.Range("N55").Value = -2482.25
# row 57
$ws.Range("H57").Value = 175000
$ws.Range("I57").Value = 300000
$ws.Range("K57").Value = 300000
$ws.Range("M57").Value = -299434
# row 61
$ws.Range("H61").Value = 7272
$ws.Range("I61").Value = 7196.2856
$ws.Range("J61").Value = 7448.6665
$ws.Range("K61").Value = 7196.2856
$ws.Range("L61").Value = 7448.6665
$ws.Range("M61").Value = -6994.2856
$ws.Range("N61").Value = -7852.6665
# row 82
$ws.Range("H82").Value = 2236.4285
$ws.Range("I82").Value = 2069.7144
$ws.Range("K82").Value = 2069.7144
$ws.Range("M82").Value = -1708.7144
# row 85
$ws.Range("H85").Value = 2236.4285
$ws.Range("I85").Value = 2069.7144
$ws.Range("K85").Value = 2069.7144
$ws.Range("M85").Value = -821.7143999999998
# row 93
$ws.Range("H93").Value = 6377.8887
$ws.Range("J93").Value = 800.6667
$ws.Range("L93").Value = 800.6667
$ws.Range("N93").Value = -3296.6667
# row 104
$ws.Range("H104").Value = 26249.75
$ws.Range("J104").Value = 26249.75
$ws.Range("L104").Value = 26249.75
$ws.Range("N104").Value = -33237.75
# row 113
$ws.Range("H113").Value = 7272
$ws.Range("I113").Value = 7196.2856
$ws.Range("J113").Value = 7448.6665
$ws.Range("K113").Value = 7196.2856
$ws.Range("L113").Value = 7448.6665
$ws.Range("M113").Value = -5026.2856
$ws.Range("N113").Value = -11788.6665

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 108
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
# row 113
$ws.Range("H113").Value = 8251.375
$ws.Range("I113").Value = 13499.75
$ws.Range("K113").Value = 40499.25
$ws.Range("M113").Value = -38329.25
